$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.420.46'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '1.568.62'
$ws.Range("E3").Value = '  -0.35%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''1.000'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Value = '''285.21'
$ws.Range("E6").Value = '  -2.31%  '
$ws.Range("D7").Value = '''0.3656'
$ws.Range("E7").Value = '  -1.98%  '
$ws.Range("D8").Value = '''48.61'
$ws.Range("E8").Value = '  -2.73%  '
$ws.Range("D9").Value = '''0.3335'
$ws.Range("E9").Value = '  -1.91%  '
$ws.Range("D10").Value = '''1.126'
$ws.Range("E10").Value = '  -1.69%  '
$ws.Range("D11").Value = '''0.07410'
$ws.Range("E11").Value = '  -2.08%  '
$ws.Range("D13").Value = '''20.77'
$ws.Range("E13").Value = '  -2.65%  '
$ws.Range("D14").Value = '''5.961'
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("D15").Value = '''6.910'
$ws.Range("E15").Value = '  -0.94%  '
$ws.Range("D16").Value = '1.569.24'
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").Value = '''0.00001105'
$ws.Range("E17").Value = '  -1.58%  '
$ws.Range("D18").Value = '''88.18'
$ws.Range("E18").Value = '  -3.05%  '
$ws.Range("D19").Value = '''0.06708'
$ws.Range("E19").Value = '  -0.62%  '
$ws.Range("D20").Value = '''1.000'
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("D21").Value = '''6.364'
$ws.Range("E21").Value = '  +0.92%  '
$ws.Range("D22").Value = '''16.20'
$ws.Range("E22").Value = '  -0.79%  '
$ws.Range("D23").Value = '''12.01'
$ws.Range("E23").Value = '  -1.38%  '
$ws.Range("D24").Value = '22.420.81'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").Value = '''2.387'
$ws.Range("E25").Value = '  +2.10%  '
$ws.Range("D26").Value = '''2.550'
$ws.Range("E26").Value = '  -5.08%  '
$ws.Range("D27").Value = '''150.00'
$ws.Range("E27").Value = '  +0.93%  '
$ws.Range("D28").Value = '''19.43'
$ws.Range("E28").Value = '  -3.41%  '
$ws.Range("D29").Value = '''5.005'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = '''124.02'
$ws.Range("E30").Value = '  -1.24%  '
$ws.Range("D31").Value = '1.744.21'
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("D32").Value = '''1.052'
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("D33").Value = '''6.116'
$ws.Range("E33").Value = '  -1.30%  '
$ws.Range("D34").Value = '''1.989'
$ws.Range("E34").Value = '  +0.18%  '
$ws.Range("D35").Value = '''9.819'
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("D36").Value = '''0.08280'
$ws.Range("E36").Value = '  -1.30%  '
$ws.Range("D37").Value = '''0.02417'
$ws.Range("D38").Value = '''0.2238'
$ws.Range("E38").Value = '  -2.75%  '
$ws.Range("D39").Value = '''0.06436'
$ws.Range("E39").Value = '  -1.29%  '
$ws.Range("D40").Value = '''5.379'
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("D41").Value = '''1.283'
$ws.Range("E41").Value = '  -6.66%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '''11.20'
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '''0.6230'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").Value = '''0.9997'
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("D45").Value = '''13.85'
$ws.Range("E45").Value = '  -1.42%  '
$ws.Range("D46").Value = '''0.6015'
$ws.Range("E46").Value = '  +3.41%  '
$ws.Range("D47").Value = '''3.752'
$ws.Range("E47").Value = '  -1.53%  '
$ws.Range("D48").Value = '''2.034'
$ws.Range("E48").Value = '  -2.07%  '
$ws.Range("D49").Value = '''123.99'
$ws.Range("E49").Value = '  -4.24%  '
$ws.Range("D50").Value = '''1.215'
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("D51").Value = '''0.07207'
$ws.Range("E51").Value = '  -1.61%  '
